$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data entered for rows 20 and 21 (columns AB, AC, AF, AI) ---
$ws.Range("AB20").Value = 532535
$ws.Range("AC20").Value = 879155
$ws.Range("AF20").Value = 390158
$ws.Range("AI20").Value = 474783
$ws.Range("AJ20").Formula = "=SUM(AB20,AC20,AF20,AI20)"

$ws.Range("AB21").Value = 739182
$ws.Range("AC21").Value = 878337
$ws.Range("AF21").Value = 462342
$ws.Range("AI21").Value = 335683
$ws.Range("AJ21").Formula = "=SUM(AB21,AC21,AF21,AI21)"

# --- Extend the AJ column formatting (right border) down through the table ---
$ws.Range("AJ15:AJ23").Borders.Item(10).Color = 0
$ws.Range("AJ15:AJ23").Borders.Item(10).Weight = 2
$ws.Range("AJ15:AJ23").Borders.Item(10).LineStyle = 1

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("AJ9").Select()
